$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 156, shifting existing rows 156-163 down to 157-164.
# Excel's row Insert carries the formatting of the row above down into the
# newly inserted row, which preserves the date style (s="2") on column D.
$ws.Rows.Item(156).Insert()

$ws.Range("A156").Value = 5
$ws.Range("B156").Value = "Macroferia Regional de Talca"
$ws.Range("C156").Value = "Maule"
$ws.Range("D156").Value = 44509
$ws.Range("E156").Value = 7
$ws.Range("F156").Value = 100112008
$ws.Range("G156").Value = "Coliflor"
$ws.Range("H156").Value = "Sin especificar"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 5000
$ws.Range("K156").Value = 600
$ws.Range("L156").Value = 600
$ws.Range("M156").Value = 600
$ws.Range("N156").Value = "`$/unidad"
$ws.Range("O156").Value = "Región del Maule"
$ws.Range("P156").Value = 600
$ws.Range("Q156").Value = 1
$ws.Range("R156").Value = "Hortaliza"
